# Updates the "cryptos" price/volume table (and one row swap) to match the
# latest scrape, mirroring the GitHub Actions data-refresh commit.
# Each entry is the target cell and its new literal text. Numeric-looking
# "Price" values are forced to Text (matching the original inline-string
# cells produced by the scraper) by stamping NumberFormat "@" before the
# write and restoring the "Normal" style afterwards, so plain-looking
# numbers like "0.999" or "6.56" don't silently become floats.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '57.380.92' },
    @{ Cell = 'E2'; Value = '  -4.74%  ' },
    @{ Cell = 'D3'; Value = '2.912.17' },
    @{ Cell = 'E3'; Value = '  -2.38%  ' },
    @{ Cell = 'D4'; Value = '0.999' },
    @{ Cell = 'E4'; Value = '  -0.12%  ' },
    @{ Cell = 'D5'; Value = '550.33' },
    @{ Cell = 'E5'; Value = '  -3.09%  ' },
    @{ Cell = 'D6'; Value = '126.74' },
    @{ Cell = 'E6'; Value = '  +1.61%  ' },
    @{ Cell = 'E7'; Value = '  +0.00%  ' },
    @{ Cell = 'E8'; Value = '  +2.72%  ' },
    @{ Cell = 'D9'; Value = '2.906.72' },
    @{ Cell = 'E9'; Value = '  -2.39%  ' },
    @{ Cell = 'D10'; Value = '0.125' },
    @{ Cell = 'E10'; Value = '  -5.61%  ' },
    @{ Cell = 'D11'; Value = '4.75' },
    @{ Cell = 'E11'; Value = '  -6.02%  ' },
    @{ Cell = 'D12'; Value = '0.436' },
    @{ Cell = 'E12'; Value = '  +0.46%  ' },
    @{ Cell = 'D13'; Value = '0.0000216' },
    @{ Cell = 'E13'; Value = '  -2.94%  ' },
    @{ Cell = 'D14'; Value = '32.22' },
    @{ Cell = 'E14'; Value = '  -0.67%  ' },
    @{ Cell = 'E15'; Value = '  +1.03%  ' },
    @{ Cell = 'D16'; Value = '3.389.44' },
    @{ Cell = 'E16'; Value = '  -2.41%  ' },
    @{ Cell = 'D17'; Value = '2.906.60' },
    @{ Cell = 'E17'; Value = '  -2.44%  ' },
    @{ Cell = 'D18'; Value = '6.56' },
    @{ Cell = 'E18'; Value = '  +6.85%  ' },
    @{ Cell = 'D19'; Value = '57.298.47' },
    @{ Cell = 'E19'; Value = '  -4.88%  ' },
    @{ Cell = 'D20'; Value = '408.46' },
    @{ Cell = 'E20'; Value = '  -4.88%  ' },
    @{ Cell = 'D21'; Value = '12.89' },
    @{ Cell = 'E21'; Value = '  -0.98%  ' },
    @{ Cell = 'D22'; Value = '0.673' },
    @{ Cell = 'E22'; Value = '  +2.25%  ' },
    @{ Cell = 'D23'; Value = '6.86' },
    @{ Cell = 'E23'; Value = '  -3.71%  ' },
    @{ Cell = 'D24'; Value = '12.72' },
    @{ Cell = 'E24'; Value = '  -1.15%  ' },
    @{ Cell = 'D25'; Value = '78.34' },
    @{ Cell = 'E25'; Value = '  -0.58%  ' },
    @{ Cell = 'E26'; Value = '  +0.22%  ' },
    @{ Cell = 'D27'; Value = '0.997' },
    @{ Cell = 'E27'; Value = '  -0.17%  ' },
    @{ Cell = 'D28'; Value = '2.46' },
    @{ Cell = 'E28'; Value = '  -1.12%  ' },
    @{ Cell = 'D29'; Value = '7.30' },
    @{ Cell = 'E29'; Value = '  +3.03%  ' },
    @{ Cell = 'D30'; Value = '1.94' },
    @{ Cell = 'E30'; Value = '  +3.38%  ' },
    @{ Cell = 'D31'; Value = '6.00' },
    @{ Cell = 'E31'; Value = '  +0.11%  ' },
    @{ Cell = 'D32'; Value = '24.79' },
    @{ Cell = 'E32'; Value = '  -1.31%  ' },
    @{ Cell = 'D33'; Value = '0.0983' },
    @{ Cell = 'E33'; Value = '  +6.45%  ' },
    @{ Cell = 'D34'; Value = '0.918' },
    @{ Cell = 'E34'; Value = '  -2.77%  ' },
    @{ Cell = 'D35'; Value = '5.48' },
    @{ Cell = 'E35'; Value = '  -0.72%  ' },
    @{ Cell = 'D36'; Value = '2.03' },
    @{ Cell = 'E36'; Value = '  -9.46%  ' },
    @{ Cell = 'D37'; Value = '48.28' },
    @{ Cell = 'E37'; Value = '  -2.03%  ' },
    @{ Cell = 'D38'; Value = '8.23' },
    @{ Cell = 'E38'; Value = '  +5.22%  ' },
    @{ Cell = 'D39'; Value = '0.0₃0640' },
    @{ Cell = 'E39'; Value = '  -0.07%  ' },
    @{ Cell = 'D40'; Value = '0.107' },
    @{ Cell = 'E40'; Value = '  +0.52%  ' },
    @{ Cell = 'D41'; Value = '0.0340' },
    @{ Cell = 'E41'; Value = '  -4.49%  ' },
    @{ Cell = 'D42'; Value = '2.46' },
    @{ Cell = 'E42'; Value = '  +4.43%  ' },
    @{ Cell = 'D43'; Value = '366.32' },
    @{ Cell = 'E43'; Value = '  -1.66%  ' },
    @{ Cell = 'D44'; Value = '2.616.44' },
    @{ Cell = 'E44'; Value = '  -0.94%  ' },
    @{ Cell = 'E45'; Value = '  -0.03%  ' },
    @{ Cell = 'D46'; Value = '120.53' },
    @{ Cell = 'E46'; Value = '  +1.45%  ' },
    @{ Cell = 'B47'; Value = 'TheGraph' },
    @{ Cell = 'C47'; Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt' },
    @{ Cell = 'D47'; Value = '0.230' },
    @{ Cell = 'E47'; Value = '  -0.98%  ' },
    @{ Cell = 'B48'; Value = 'Stellar' },
    @{ Cell = 'C48'; Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm' },
    @{ Cell = 'D48'; Value = '0.108' },
    @{ Cell = 'E48'; Value = '  +2.21%  ' },
    @{ Cell = 'D49'; Value = '1.95' },
    @{ Cell = 'E49'; Value = '  +0.39%  ' },
    @{ Cell = 'D50'; Value = '22.70' },
    @{ Cell = 'E50'; Value = '  -2.32%  ' },
    @{ Cell = 'D51'; Value = '1.95' },
    @{ Cell = 'E51'; Value = '  -1.00%  ' }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.Value -match '^[0-9]+(\.[0-9]+)?$') {
        # Purely numeric text (e.g. "0.999", "6.56") - keep it a string.
        $rng.NumberFormat = '@'
        $rng.Value = $u.Value
        $rng.Style = 'Normal'
    } else {
        $rng.Value = $u.Value
    }
}
